# Apply the "Kode Fixed Income" update (OBL00107 -> OBL00108) to both
# DPLKINV108-001 and DPLKINV108-002 sheets, and flip the active sheet /
# selection from sheet 1 to sheet 2, matching the author's edit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKINV108-001")
$ws2 = $wb.Worksheets.Item("DPLKINV108-002")

# --- Sheet 1: DPLKINV108-001 ---
$ws1.Range("M2").Value = "OBL00108"
$ws1.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nKode Fixed Income : OBL00108;`nStatus Verifikasi : 1 : Setuju`nKeterangan Verifikasi : DATA APPROVAL"

# --- Sheet 2: DPLKINV108-002 ---
$ws2.Range("M2").Value = "OBL00108"
$ws2.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nKode Fixed Income : OBL00108;`nStatus Verifikasi : 0 : Kembalikan ke Data Entry;`nKeterangan Verifikasi : DATA AKAN DIKEMBALIKAN UNTUK DIEDIT"

# --- View state: move the active tab / selection from sheet1 to sheet2 ---
$ws1.Select()
$ws1.Application.ActiveWindow.ScrollColumn = 4
$ws1.Range("G2").Select()

$ws2.Select()
$ws2.Application.ActiveWindow.ScrollColumn = 5
$ws2.Range("N2").Select()

$wb.Save()
